$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "B123"
$ws.Range("B3").Value = "kk"

$ws.Range("B3").Select()
